$d = $word.ActiveDocument
$tbl = $d.Tables(9)

$cells = @(
    $tbl.Cell(2,2),
    $tbl.Cell(2,3),
    $tbl.Cell(3,2),
    $tbl.Cell(3,3)
)

foreach ($c in $cells) {
    $r = $c.Range
    $start = $r.Start
    $end = $r.End
    # Exclude the trailing cell/paragraph mark from the range so that
    # deleting it removes the run entirely instead of leaving an empty one.
    $narrow = $d.Range($start, $end - 1)
    if ($narrow.Text -ne "") {
        $narrow.Delete()
    }
}
